$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.041699186451632
$ws.Cells.Item(2, 4).Value = 1.042967135455924
$ws.Cells.Item(2, 5).Value = 1.039923283595741
$ws.Cells.Item(2, 9).Value = 1.041867214018811
$ws.Cells.Item(2, 10).Value = 1.046778987644194
$ws.Cells.Item(2, 11).Value = 1.045742329522629
$ws.Cells.Item(2, 12).Value = 1.042707093143562
$ws.Cells.Item(2, 14).Value = 1.048265533873876

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.042778172017404
$ws.Cells.Item(3, 4).Value = 1.043775313598928
$ws.Cells.Item(3, 5).Value = 1.040844816310789
$ws.Cells.Item(3, 9).Value = 1.04219372257804
$ws.Cells.Item(3, 10).Value = 1.047503464140105
$ws.Cells.Item(3, 11).Value = 1.04636145411738
$ws.Cells.Item(3, 12).Value = 1.043438650924582
$ws.Cells.Item(3, 14).Value = 1.04899103920951

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.043476121564872
$ws.Cells.Item(4, 4).Value = 1.044297989123603
$ws.Cells.Item(4, 5).Value = 1.041441250279764
$ws.Cells.Item(4, 9).Value = 1.042403514909571
$ws.Cells.Item(4, 10).Value = 1.04797145537686
$ws.Cells.Item(4, 11).Value = 1.046761148605135
$ws.Cells.Item(4, 12).Value = 1.043911532047782
$ws.Cells.Item(4, 14).Value = 1.049459695047501

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.043769485877176
$ws.Cells.Item(5, 4).Value = 1.04451765699492
$ws.Cells.Item(5, 5).Value = 1.041692024596478
$ws.Cells.Item(5, 9).Value = 1.042491357104337
$ws.Cells.Item(5, 10).Value = 1.048168008921229
$ws.Cells.Item(5, 11).Value = 1.04692895967638
$ws.Cells.Item(5, 12).Value = 1.044110214918594
$ws.Cells.Item(5, 14).Value = 1.049656527720449

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.043818739902678
$ws.Cells.Item(6, 4).Value = 1.044554536380418
$ws.Cells.Item(6, 5).Value = 1.041734132664835
$ws.Cells.Item(6, 9).Value = 1.042506085405549
$ws.Cells.Item(6, 10).Value = 1.04820099998598
$ws.Cells.Item(6, 11).Value = 1.046957122960517
$ws.Cells.Item(6, 12).Value = 1.044143567821457
$ws.Cells.Item(6, 14).Value = 1.049689565636296

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.04348004172347
$ws.Cells.Item(7, 4).Value = 1.044300924591967
$ws.Cells.Item(7, 5).Value = 1.041444601007178
$ws.Cells.Item(7, 9).Value = 1.042404690054299
$ws.Cells.Item(7, 10).Value = 1.047974082481411
$ws.Cells.Item(7, 11).Value = 1.046763391771415
$ws.Cells.Item(7, 12).Value = 1.043914187315321
$ws.Cells.Item(7, 14).Value = 1.049462325882842

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.042063881820355
$ws.Cells.Item(8, 4).Value = 1.043240318702295
$ws.Cells.Item(8, 5).Value = 1.040234690641353
$ws.Cells.Item(8, 9).Value = 1.041977865636684
$ws.Cells.Item(8, 10).Value = 1.047023992164989
$ws.Cells.Item(8, 11).Value = 1.045951756110284
$ws.Cells.Item(8, 12).Value = 1.042954427036203
$ws.Cells.Item(8, 14).Value = 1.048510886329194

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.039566671718287
$ws.Cells.Item(9, 4).Value = 1.041369343616321
$ws.Cells.Item(9, 5).Value = 1.038103759798202
$ws.Cells.Item(9, 9).Value = 1.041214409476974
$ws.Cells.Item(9, 10).Value = 1.04534373400741
$ws.Cells.Item(9, 11).Value = 1.044514504818132
$ws.Cells.Item(9, 12).Value = 1.041259490078926
$ws.Cells.Item(9, 14).Value = 1.046828242012303

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.037900650148182
$ws.Cells.Item(10, 4).Value = 1.040120663708747
$ws.Cells.Item(10, 5).Value = 1.036683877665217
$ws.Cells.Item(10, 9).Value = 1.040697817559755
$ws.Cells.Item(10, 10).Value = 1.044219472301486
$ws.Cells.Item(10, 11).Value = 1.043551603954343
$ws.Cells.Item(10, 12).Value = 1.040127034603456
$ws.Cells.Item(10, 14).Value = 1.045702383725802

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.037178947593417
$ws.Cells.Item(11, 4).Value = 1.039579650819083
$ws.Cells.Item(11, 5).Value = 1.036069229013195
$ws.Cells.Item(11, 9).Value = 1.040472319613455
$ws.Cells.Item(11, 10).Value = 1.043731682934818
$ws.Cells.Item(11, 11).Value = 1.043133534544205
$ws.Cells.Item(11, 12).Value = 1.039636075944033
$ws.Cells.Item(11, 14).Value = 1.045213901642283

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.036910828559326
$ws.Cells.Item(12, 4).Value = 1.039378645638569
$ws.Cells.Item(12, 5).Value = 1.03584094658967
$ws.Cells.Item(12, 9).Value = 1.040388287560958
$ws.Cells.Item(12, 10).Value = 1.04355034926831
$ws.Cells.Item(12, 11).Value = 1.042978075529497
$ws.Cells.Item(12, 12).Value = 1.039453622049256
$ws.Cells.Item(12, 14).Value = 1.045032310461169

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.036968343123664
$ws.Cells.Item(13, 4).Value = 1.039421764143429
$ws.Cells.Item(13, 5).Value = 1.035889912793306
$ws.Cells.Item(13, 9).Value = 1.04040632503517
$ws.Cells.Item(13, 10).Value = 1.043589252610376
$ws.Cells.Item(13, 11).Value = 1.043011429684253
$ws.Cells.Item(13, 12).Value = 1.039492763103674
$ws.Cells.Item(13, 14).Value = 1.045071269050443

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.037156785745474
$ws.Cells.Item(14, 4).Value = 1.039563036665175
$ws.Cells.Item(14, 5).Value = 1.036050358597859
$ws.Cells.Item(14, 9).Value = 1.040465379048652
$ws.Cells.Item(14, 10).Value = 1.043716696837114
$ws.Cells.Item(14, 11).Value = 1.043120687710313
$ws.Cells.Item(14, 12).Value = 1.039620996084801
$ws.Cells.Item(14, 14).Value = 1.045198894262601

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.03727288526708
$ws.Cells.Item(15, 4).Value = 1.039650072851296
$ws.Cells.Item(15, 5).Value = 1.036149217924967
$ws.Cells.Item(15, 9).Value = 1.040501728119594
$ws.Cells.Item(15, 10).Value = 1.043795199948029
$ws.Cells.Item(15, 11).Value = 1.043187982732996
$ws.Cells.Item(15, 12).Value = 1.039699992727754
$ws.Cells.Item(15, 14).Value = 1.045277508856938

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.037948540686426
$ws.Cells.Item(16, 4).Value = 1.040156562088287
$ws.Cells.Item(16, 5).Value = 1.036724673443635
$ws.Cells.Item(16, 9).Value = 1.040712744967038
$ws.Cells.Item(16, 10).Value = 1.044251824677572
$ws.Cells.Item(16, 11).Value = 1.043579326091105
$ws.Cells.Item(16, 12).Value = 1.040159605301453
$ws.Cells.Item(16, 14).Value = 1.045734782045973

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.038372279421554
$ws.Cells.Item(17, 4).Value = 1.040474182138774
$ws.Cells.Item(17, 5).Value = 1.037085686888946
$ws.Cells.Item(17, 9).Value = 1.040844625533283
$ws.Cells.Item(17, 10).Value = 1.044537991660217
$ws.Cells.Item(17, 11).Value = 1.043824503642427
$ws.Cells.Item(17, 12).Value = 1.04044774792807
$ws.Cells.Item(17, 14).Value = 1.046021355418556

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.038619409761947
$ws.Cells.Item(18, 4).Value = 1.040659413055007
$ws.Cells.Item(18, 5).Value = 1.037296276298392
$ws.Cells.Item(18, 9).Value = 1.040921374581952
$ws.Cells.Item(18, 10).Value = 1.04470481380599
$ws.Cells.Item(18, 11).Value = 1.043967402829832
$ws.Cells.Item(18, 12).Value = 1.040615758962489
$ws.Cells.Item(18, 14).Value = 1.046188414470911

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.03870366987901
$ws.Cells.Item(19, 4).Value = 1.040722566666718
$ws.Cells.Item(19, 5).Value = 1.037368084623468
$ws.Cells.Item(19, 9).Value = 1.040947514403009
$ws.Cells.Item(19, 10).Value = 1.044761679873903
$ws.Cells.Item(19, 11).Value = 1.044016109320508
$ws.Cells.Item(19, 12).Value = 1.040673036603661
$ws.Cells.Item(19, 14).Value = 1.046245361295163

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.03832681930165
$ws.Cells.Item(20, 4).Value = 1.040440107769059
$ws.Cells.Item(20, 5).Value = 1.037046951869313
$ws.Cells.Item(20, 9).Value = 1.040830494062649
$ws.Cells.Item(20, 10).Value = 1.044507298397866
$ws.Cells.Item(20, 11).Value = 1.043798209658628
$ws.Cells.Item(20, 12).Value = 1.04041683891876
$ws.Cells.Item(20, 14).Value = 1.045990618568253

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.037101295382718
$ws.Cells.Item(21, 4).Value = 1.039521436777787
$ws.Cells.Item(21, 5).Value = 1.036003110605129
$ws.Cells.Item(21, 9).Value = 1.04044799662503
$ws.Cells.Item(21, 10).Value = 1.0436791717481
$ws.Cells.Item(21, 11).Value = 1.043088518622534
$ws.Cells.Item(21, 12).Value = 1.039583237153337
$ws.Cells.Item(21, 14).Value = 1.045161315883657

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.036330489976467
$ws.Cells.Item(22, 4).Value = 1.03894354896784
$ws.Cells.Item(22, 5).Value = 1.035346953414138
$ws.Cells.Item(22, 9).Value = 1.040205931109257
$ws.Cells.Item(22, 10).Value = 1.043157645166809
$ws.Cells.Item(22, 11).Value = 1.04264132742765
$ws.Cells.Item(22, 12).Value = 1.039058597914271
$ws.Cells.Item(22, 14).Value = 1.044639048674801

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.036739134377262
$ws.Cells.Item(23, 4).Value = 1.039249925009395
$ws.Cells.Item(23, 5).Value = 1.035694780813301
$ws.Cells.Item(23, 9).Value = 1.040334403870848
$ws.Cells.Item(23, 10).Value = 1.043434196975663
$ws.Cells.Item(23, 11).Value = 1.042878484895563
$ws.Cells.Item(23, 12).Value = 1.039336768511159
$ws.Cells.Item(23, 14).Value = 1.044915993218944

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.038347360863226
$ws.Cells.Item(24, 4).Value = 1.040455504607192
$ws.Cells.Item(24, 5).Value = 1.037064454506619
$ws.Cells.Item(24, 9).Value = 1.040836880005495
$ws.Cells.Item(24, 10).Value = 1.044521167652876
$ws.Cells.Item(24, 11).Value = 1.043810091113916
$ws.Cells.Item(24, 12).Value = 1.040430805547691
$ws.Cells.Item(24, 14).Value = 1.046004507519196

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.040212471967566
$ws.Cells.Item(25, 4).Value = 1.041853276980113
$ws.Cells.Item(25, 5).Value = 1.038654526946066
$ws.Cells.Item(25, 9).Value = 1.041413124555788
$ws.Cells.Item(25, 10).Value = 1.04577884157871
$ws.Cells.Item(25, 11).Value = 1.044886903268138
$ws.Cells.Item(25, 12).Value = 1.041698112203022
$ws.Cells.Item(25, 14).Value = 1.047263967486261
